# Append the latest daily profit record (run on 2025-11-10) as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 85

# Column A holds the date as literal text (matching the existing rows),
# e.g. "11/10/2025". Force the cell to Text first so Excel doesn't
# auto-convert the date-looking string into a date serial number, then
# clear the formatting afterwards so the cell ends up with the same
# "no explicit style" look as the rest of the date column.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "11/10/2025"
$ws.Range("A" + $newRow).ClearFormats()

# Column B holds the plain numeric profit value.
$ws.Range("B" + $newRow).Value = 10671.17
